$d = $word.ActiveDocument

$d.Content.Find.Execute("(10 credits: ~300 words)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(10 credits: ~350 words)", 2)

$d.Content.Find.Execute("(15 credits: ~400 words)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(15 credits: ~475 words)", 2)

$d.Content.Find.Execute("(20 credits: ~500 words)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(20 credits: ~600 words)", 2)
